$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead text: bump the report Volume/Number and the week-of dates ---
# Only the trailing run of each rich-text cell changes; do the rightmost date run
# first so its start offset is not shifted by the length change of the earlier run.
$ws.Range("A8").Characters(21, 2).Text = "29"
$ws.Range("C9").Characters(46, 9).Text = "7/21/2024"
$ws.Range("C9").Characters(27, 8).Text = "7/15/2024"

# --- Crime-grid updates (rows 14-30): refreshed weekly figures ---

# A handful of cells flip from a number to the sheets "not applicable" text markers
# ("0" / "***.*"). Copying formats first from a stable cell that already carries that
# text style keeps the destination on the existing text style (style 14) instead of
# Excel minting a new quote-prefixed numeric style; copying values second then carries
# over the literal shared string instead of a number.
$zeroSrc = $ws.Range("D14")
$naSrc = $ws.Range("N14")

$zeroSrc.Copy()
$ws.Range("C14").PasteSpecial(-4122)
$zeroSrc.Copy()
$ws.Range("C14").PasteSpecial(-4163)
$zeroSrc.Copy()
$ws.Range("D15").PasteSpecial(-4122)
$zeroSrc.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$naSrc.Copy()
$ws.Range("E15").PasteSpecial(-4122)
$naSrc.Copy()
$ws.Range("E15").PasteSpecial(-4163)
$zeroSrc.Copy()
$ws.Range("C23").PasteSpecial(-4122)
$zeroSrc.Copy()
$ws.Range("C23").PasteSpecial(-4163)
$zeroSrc.Copy()
$ws.Range("D23").PasteSpecial(-4122)
$zeroSrc.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$naSrc.Copy()
$ws.Range("E23").PasteSpecial(-4122)
$naSrc.Copy()
$ws.Range("E23").PasteSpecial(-4163)
$zeroSrc.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$zeroSrc.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$naSrc.Copy()
$ws.Range("E27").PasteSpecial(-4122)
$naSrc.Copy()
$ws.Range("E27").PasteSpecial(-4163)
$zeroSrc.Copy()
$ws.Range("C28").PasteSpecial(-4122)
$zeroSrc.Copy()
$ws.Range("C28").PasteSpecial(-4163)
$zeroSrc.Copy()
$ws.Range("D28").PasteSpecial(-4122)
$zeroSrc.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$naSrc.Copy()
$ws.Range("E28").PasteSpecial(-4122)
$naSrc.Copy()
$ws.Range("E28").PasteSpecial(-4163)
$zeroSrc.Copy()
$ws.Range("G29").PasteSpecial(-4122)
$zeroSrc.Copy()
$ws.Range("G29").PasteSpecial(-4163)
$naSrc.Copy()
$ws.Range("H29").PasteSpecial(-4122)
$naSrc.Copy()
$ws.Range("H29").PasteSpecial(-4163)
$zeroSrc.Copy()
$ws.Range("G30").PasteSpecial(-4122)
$zeroSrc.Copy()
$ws.Range("G30").PasteSpecial(-4163)
$naSrc.Copy()
$ws.Range("H30").PasteSpecial(-4122)
$naSrc.Copy()
$ws.Range("H30").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Plain numeric updates
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -47.058823529411
$ws.Range("I16").Value = 95
$ws.Range("J16").Value = 83
$ws.Range("K16").Value = 14.457831325301
$ws.Range("L16").Value = -11.214953271028
$ws.Range("M16").Value = -29.629629629629
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -55.555555555555
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 141
$ws.Range("J17").Value = 138
$ws.Range("K17").Value = 2.173913043478
$ws.Range("L17").Value = -14.024390243902
$ws.Range("M17").Value = 45.360824742268
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 4
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 55
$ws.Range("J18").Value = 73
$ws.Range("K18").Value = -24.657534246575
$ws.Range("L18").Value = -3.508771929824
$ws.Range("M18").Value = 1.851851851851
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -15.625
$ws.Range("I19").Value = 194
$ws.Range("J19").Value = 181
$ws.Range("K19").Value = 7.182320441988
$ws.Range("L19").Value = 15.476190476190
$ws.Range("M19").Value = 31.972789115646
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 7
$ws.Range("H20").Value = -56.25
$ws.Range("J20").Value = 83
$ws.Range("K20").Value = -28.915662650602
$ws.Range("L20").Value = -28.915662650602
$ws.Range("M20").Value = 51.282051282051
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -46.153846153846
$ws.Range("F21").Value = 68
$ws.Range("G21").Value = 98
$ws.Range("H21").Value = -30.612244897959
$ws.Range("I21").Value = 554
$ws.Range("J21").Value = 567
$ws.Range("K21").Value = -2.292768959435
$ws.Range("L21").Value = -5.621805792163
$ws.Range("M21").Value = 14.226804123711
$ws.Range("F22").Value = 1
$ws.Range("L22").Value = -61.111111111111
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 200
$ws.Range("L23").Value = 13.333333333333
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -35.294117647058
$ws.Range("F24").Value = 54
$ws.Range("G24").Value = 69
$ws.Range("H24").Value = -21.739130434782
$ws.Range("I24").Value = 491
$ws.Range("J24").Value = 525
$ws.Range("K24").Value = -6.476190476190
$ws.Range("L24").Value = -41.547619047619
$ws.Range("M24").Value = 57.371794871794
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = -15.789473684210
$ws.Range("I25").Value = 147
$ws.Range("J25").Value = 134
$ws.Range("K25").Value = 9.701492537313
$ws.Range("L25").Value = -72.106261859582
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = -6.451612903225
$ws.Range("I26").Value = 191
$ws.Range("J26").Value = 222
$ws.Range("K26").Value = -13.963963963964
$ws.Range("L26").Value = -18.025751072961
$ws.Range("M26").Value = -28.464419475655
